# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for handoff"
# and refreshes the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
# timestamps, then re-fits the Status columns to the new (longer) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1. Status: "In Translation" -> "Ready for handoff" everywhere it is shown.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# 2. Refresh the generate/handoff timestamps for this handback run.
$wsZhCn.Range("H2").Value = "2016-09-01 00:42:36"
$wsOverview.Range("G2").Value = "2016-09-01 00:42:41"
$wsDeDe.Range("H2").Value = "2016-09-01 00:42:41"

# 3. Re-fit the Status columns now that "Ready for handoff" is longer than
#    "In Translation".
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
